$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name reshuffle (mirrors sharedStrings.xml <si> reordering) ---
$ws.Range("A137").Value = 'Uganda'
$ws.Range("A138").Value = 'Crucero'
$ws.Range("A146").Value = 'Benin'
$ws.Range("A147").Value = 'Togo'
$ws.Range("A148").Value = 'Tanzania'
$ws.Range("A149").Value = 'Suazilandia'
$ws.Range("A150").Value = 'Estado de Palestina'
$ws.Range("A151").Value = 'Liberia'
$ws.Range("A152").Value = 'Reunion'
$ws.Range("A190").Value = 'Namibia'
$ws.Range("A191").Value = 'Guam'

# --- Updated case statistics ---
# Row 4
$ws.Range("B4").Value = 2182951
$ws.Range("C4").Value = 1
$ws.Range("E4").Value = 1174802
# Row 21
$ws.Range("B21").Value = 94481
$ws.Range("C21").Value = 3862
$ws.Range("D21").Value = 36264
$ws.Range("E21").Value = 56955
$ws.Range("G21").Value = 53
$ws.Range("H21").Value = 1262
# Row 25
$ws.Range("B25").Value = 60155
$ws.Range("C25").Value = 55
$ws.Range("D25").Value = 16625
$ws.Range("E25").Value = 33867
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 9663
# Row 34
$ws.Range("B34").Value = 40400
$ws.Range("C34").Value = 1106
$ws.Range("D34").Value = 15703
$ws.Range("E34").Value = 22466
$ws.Range("G34").Value = 33
$ws.Range("H34").Value = 2231
# Row 41
$ws.Range("B41").Value = 26781
$ws.Range("C41").Value = 361
$ws.Range("D41").Value = 6552
$ws.Range("E41").Value = 19126
$ws.Range("G41").Value = 5
$ws.Range("H41").Value = 1103
# Row 42
$ws.Range("B42").Value = 25623
$ws.Range("C42").Value = 96
$ws.Range("D42").Value = 5506
$ws.Range("E42").Value = 19626
$ws.Range("G42").Value = 13
$ws.Range("H42").Value = 491
# Row 51
$ws.Range("E51").Value = 5699
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 47
# Row 68
$ws.Range("B68").Value = 8921
$ws.Range("C68").Value = 36
$ws.Range("D68").Value = 7880
$ws.Range("E68").Value = 829
# Row 70
$ws.Range("B70").Value = 8505
$ws.Range("C70").Value = 11
$ws.Range("D70").Value = 7733
$ws.Range("E70").Value = 651
# Row 73
$ws.Range("B73").Value = 7112
$ws.Range("C73").Value = 4
$ws.Range("E73").Value = 586
# Row 103
$ws.Range("D103").Value = 1371
$ws.Range("E103").Value = 523
# Row 107
$ws.Range("B107").Value = 1776
$ws.Range("C107").Value = 3
$ws.Range("D107").Value = 1441
$ws.Range("E107").Value = 259
# Row 110
$ws.Range("B110").Value = 1672
$ws.Range("C110").Value = 82
$ws.Range("D110").Value = 1064
$ws.Range("E110").Value = 571
$ws.Range("G110").Value = 1
$ws.Range("H110").Value = 37
# Row 113
$ws.Range("B113").Value = 1499
$ws.Range("C113").Value = 3
$ws.Range("E113").Value = 31
# Row 123
$ws.Range("D123").Value = 1069
$ws.Range("E123").Value = 40
# Row 137
$ws.Range("B137").Value = 724
$ws.Range("C137").Value = 19
$ws.Range("D137").Value = 351
$ws.Range("E137").Value = 373
$ws.Range("H137").Value = 0
# Row 138
$ws.Range("B138").Value = 712
$ws.Range("D138").Value = 651
$ws.Range("E138").Value = 48
$ws.Range("H138").Value = 13
# Row 146
$ws.Range("C146").Value = 49
$ws.Range("D146").Value = 236
$ws.Range("E146").Value = 287
$ws.Range("H146").Value = 9
# Row 147
$ws.Range("B147").Value = 532
$ws.Range("C147").Value = 1
$ws.Range("D147").Value = 308
$ws.Range("E147").Value = 211
$ws.Range("H147").Value = 13
# Row 148
$ws.Range("B148").Value = 509
$ws.Range("D148").Value = 183
$ws.Range("E148").Value = 305
$ws.Range("H148").Value = 21
# Row 149
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 249
$ws.Range("E149").Value = 253
$ws.Range("H149").Value = 4
# Row 150
$ws.Range("B150").Value = 506
$ws.Range("C150").Value = 1
$ws.Range("D150").Value = 415
$ws.Range("E150").Value = 88
$ws.Range("H150").Value = 3
# Row 151
$ws.Range("B151").Value = 498
$ws.Range("D151").Value = 221
$ws.Range("E151").Value = 244
$ws.Range("H151").Value = 33
# Row 152
$ws.Range("B152").Value = 496
$ws.Range("D152").Value = 460
$ws.Range("E152").Value = 35
$ws.Range("H152").Value = 1
# Row 190
$ws.Range("B190").Value = 34
$ws.Range("C190").Value = 2
$ws.Range("D190").Value = 18
$ws.Range("E190").Value = 16
$ws.Range("H190").Value = 0
# Row 191
$ws.Range("D191").Value = 0
$ws.Range("E191").Value = 31
$ws.Range("H191").Value = 1

# --- Timestamp update ---
$ws.Range("A1").Value = 'Datos actualizados a 16 de Junio de 2020 a las 12:01'
